# Apply the changes described by the commit:
# "added the instructions for testing CSV upload and report generation.
#  Modified the sprint planning and burndown charts to match the commits"

$wb  = $excel.ActiveWorkbook
$wsActual   = $wb.Worksheets.Item("Actual")
$wsBurndown = $wb.Worksheets.Item("Burndown")

# ---------------------------------------------------------------------------
# "Actual" sheet (Table1) - rows 2 & 3: the "V:6" note was entered under the
# wrong table column ("5" / column I) - move it one column to the right, into
# the "6" column (column J), where it belongs.
# ---------------------------------------------------------------------------
$wsActual.Range("J2").Value = $wsActual.Range("I2").Value2
$wsActual.Range("I2").ClearContents()

$wsActual.Range("J3").Value = $wsActual.Range("I3").Value2
$wsActual.Range("I3").ClearContents()

# Row 7 - Dependency column (C7) was left empty; record task "T24" there.
$wsActual.Range("C7").Value = "T24"

# ---------------------------------------------------------------------------
# "Burndown" sheet - the "Actual (in story points)" row had a stale value
# (11) for day 5; update it to match the real burndown (23).
# ---------------------------------------------------------------------------
$wsBurndown.Range("G5").Value = 23

# ---------------------------------------------------------------------------
# Update which sheet/cell is active & selected in each sheet view, matching
# what was left selected after the edits above: "Actual" becomes the active
# (displayed) tab with D10 selected, and "Burndown" is no longer the active
# tab, with G5 left selected.
# ---------------------------------------------------------------------------
$wsBurndown.Range("G5").Select() | Out-Null
$wsActual.Activate()
$wsActual.Range("D10").Select() | Out-Null
